$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("U32_01", 32, "Hình dáng", "Shape", "The ruler has a rectangular shape", "a rectangular shape / hình chữ nhật", "N"),
  @("U32_02", 32, "Cổ áo", "Collar", "My shirt collar is stained / … bị ố", "a shirt collar", "N"),
  @("U32_03", 32, "Đường cong", "Curve", "The fast car took a sharp curve", "a sharp curve / khúc cua gấp", "N"),
  @("U32_04", 32, "Khuy áo", "Button", "One of your shirt buttons is missing", "a shirt button / một cái khuy áo", "N"),
  @("U32_05", 32, "Túi", "Pocket", "The keys are in my pocket", "in one's pocket", "N"),
  @("U32_06", 32, "Tay áo", "Sleeve", "His shirt has short sleeves", "a short sleeve / áo cộc tay", "N"),
  @("U32_07", 32, "Tủ quần áo", "Closet", "the child is hiding in the closet", "in the closet", "N"),
  @("U32_08", 32, "Tất", "Sock", "Your christmas gift is a pair of socks", "a pair of socks", "N"),
  @("U32_09", 32, "Trang sức", "Juwelry", "Gold jewelry is expensive", "gold jewelry", "N"),
  @("U32_10", 32, "Vòng tay", "Bracelet", "she wears a diamond bracelet on her wrist", "a diamond bracelet", "N"),
  @("U32_11", 32, "Cửa hàng", "Store", "Apple now has retail stores in vietnam", "a retail store / cửa hàng bán lẻ", "N"),
  @("U32_12", 32, "Thắt", "Tie", "I don’t know how to tie a tie", "tie a tie / thắt cà vạt", "V"),
  @("U32_13", 32, "Da", "Leather", "Bikers usually wear leather jackets", "a leather jacket / áo khoác da", "N"),
  @("U32_14", 32, "Vẻ đẹp", "Beauty", "The beauty of this view amazes me", "beauty of something or somebody", "N"),
  @("U32_15", 32, "Hợp", "Suit", "This hat suits you well", "suit somebody well", "V"),
  @("U32_16", 32, "Thử đồ", "Try on", "Why don’t you try on the dress?", "try on something", "V"),
  @("U32_17", 32, "Lỏng", "Loose", "can you fix my shirt's loose button?", "a loose button", "Adj"),
  @("U32_18", 32, "Quần", "Pants", "You need a new pair of pants", "a pair of pants / một chiếc quần dài", "N"),
  @("U32_19", 32, "Chật", "Tight", "These pants are too tight for me", "too tight for somebody", "Adj"),
  @("U32_20", 32, "Phần eo", "Waist", "His clothes are wet from the waist down", "from the waist down / từ thắt lưng trở xuống", "N"),
  @("U32_21", 32, "Kim cương", "Diamond", "The diamond rind he gave her was huge", "a diamond ring", "N"),
  @("U32_22", 32, "Mẫu họa tiết", "Pattern", "I like the pattern on your shirt", "the pattern on something", "N"),
  @("U32_23", 32, "Hợp với", "Go with", "My new shirt goes with my favorite watch", "go with something", "V"),
  @("U32_24", 32, "Mảnh", "Piece", "This purse can go with any piece of clothing", "a piece of clothing / món đồ thời trang", "N"),
  @("U32_25", 32, "Khoe mẽ", "Show off", "He is showing off his new watch", "show off something", "V"),
  @("U32_26", 32, "Ngầu", "Cool", "The new poster design looks cool", "to look cool / trong thật ngầu", "Adj"),
  @("U32_27", 32, "Nhiều màu sắc", "Colorful", "Children always make the most colorful paintings", "A colorful painting / một bức tranh sặc sỡ", "Adj"),
  @("U32_28", 32, "Thời trang", "Fashion", "She writes about the latest fashion trends", "a fashion trend", "N"),
  @("U32_29", 32, "Thanh lịch", "Chic", "The fashionista has a chic style", "a chic style / phong cách thanh lịch", "Adj"),
  @("U32_30", 32, "Gói lại", "Wrap up", "The cashier wrapped up the shirt for me", "wrap up something", "V"),
)

$startRow = 932
for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $startRow + $i
  $row = $data[$i]
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Cells.Item($r, 6).Value = $row[5]
  $ws.Cells.Item($r, 7).Value = $row[6]
}

$ws.Range("B930").Select()
$ws.Range("C953").Select()

Write-Host "Done writing unit 32 rows"